$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "E" (error) column values for existing rows 3..21 — the
# table now reflects a halving error term starting at 2.5 instead of 1.0.
# Values are entered with a leading apostrophe so they stay text cells
# (matching the rest of this inline-string table) instead of being
# auto-converted to numeric cells.
$ws.Range("D3").Value  = "'2.5"
$ws.Range("D4").Value  = "'1.25"
$ws.Range("D5").Value  = "'0.625"
$ws.Range("D6").Value  = "'0.3125"
$ws.Range("D7").Value  = "'0.15625"
$ws.Range("D8").Value  = "'0.078125"
$ws.Range("D9").Value  = "'0.0390625"
$ws.Range("D10").Value = "'0.01953125"
$ws.Range("D11").Value = "'0.009765625"
$ws.Range("D12").Value = "'0.0048828125"
$ws.Range("D13").Value = "'0.00244140625"
$ws.Range("D14").Value = "'0.001220703125"
$ws.Range("D15").Value = "'0.0006103515625"
$ws.Range("D16").Value = "'0.00030517578125"
$ws.Range("D17").Value = "'0.000152587890625"
$ws.Range("D18").Value = "'7.62939453125e-05"
$ws.Range("D19").Value = "'3.814697265625e-05"
$ws.Range("D20").Value = "'1.9073486328125e-05"
$ws.Range("D21").Value = "'9.5367431640625e-06"

# Append the new iteration row (row 22) with the next bisection step.
$ws.Range("A22").Value = "'20"
$ws.Range("B22").Value = "'-3.1642484664917"
$ws.Range("C22").Value = "'2.51429654994448e-05"
$ws.Range("D22").Value = "'4.76837158203125e-06"

# The leading apostrophe above marks the style as quote-prefixed text;
# restore the plain default style so these cells match the rest of the
# (unstyled) data rows instead of picking up a text-quote style.
$ws.Range("D3:D21").Style = "Normal"
$ws.Range("A22:D22").Style = "Normal"
